# Auto-generated edit script applying cached-value updates described by the diff.
# The workbook stores plain numeric cache values (no formulas) in columns H-N of each
# leve-profit sheet; this script rewrites each changed cell to its new value, and clears
# the handful of cells that the diff removes entirely (value -> blank).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 88.27273
$ws.Range("I6").Value = 53.22222
$ws.Range("K6").Value = 159.66666
$ws.Range("M6").Value = -47.66666000000001
$ws.Range("H8").Value = 37
$ws.Range("I8").Value = 38.909092
$ws.Range("K8").Value = 116.727276
$ws.Range("M8").Value = 22.272724
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H38").Value = 29
$ws.Range("I38").Value = 29
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 87
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 285
$ws.Range("N38").ClearContents()
$ws.Range("H43").Value = 2390.75
$ws.Range("I43").Value = 1380.909
$ws.Range("J43").Value = 13499
$ws.Range("K43").Value = 1380.909
$ws.Range("L43").Value = 13499
$ws.Range("N43").Value = -13637
$ws.Range("M43").Value = -1311.909
$ws.Range("H62").Value = 4575.1665
$ws.Range("I62").Value = 3733.5557
$ws.Range("K62").Value = 3733.5557
$ws.Range("M62").Value = -3109.5557
$ws.Range("H65").Value = 4575.1665
$ws.Range("I65").Value = 3733.5557
$ws.Range("K65").Value = 18667.7785
$ws.Range("M65").Value = -15547.7785
$ws.Range("H70").Value = 1520.4375
$ws.Range("I70").Value = 841.8570999999999
$ws.Range("J70").Value = 2048.2222
$ws.Range("K70").Value = 2525.5713
$ws.Range("L70").Value = 6144.6666
$ws.Range("M70").Value = -2255.5713
$ws.Range("N70").Value = -6684.6666
$ws.Range("H73").Value = 1520.4375
$ws.Range("I73").Value = 841.8570999999999
$ws.Range("J73").Value = 2048.2222
$ws.Range("K73").Value = 2525.5713
$ws.Range("L73").Value = 6144.6666
$ws.Range("M73").Value = -1589.5713
$ws.Range("N73").Value = -8016.6666
$ws.Range("H80").Value = 3193.4285
$ws.Range("I80").Value = 859.6
$ws.Range("K80").Value = 2578.8
$ws.Range("M80").Value = -1580.8
$ws.Range("H83").Value = 3193.4285
$ws.Range("I83").Value = 859.6
$ws.Range("K83").Value = 7736.400000000001
$ws.Range("M83").Value = -2744.400000000001
$ws.Range("H103").Value = 430
$ws.Range("I103").Value = 504.2857
$ws.Range("K103").Value = 1512.8571
$ws.Range("M103").Value = -926.8571000000002
$ws.Range("H111").Value = 610.1667
$ws.Range("I111").Value = 676.3333
$ws.Range("J111").Value = 544
$ws.Range("K111").Value = 2028.9999
$ws.Range("L111").Value = 1632
$ws.Range("M111").Value = 1038.0001
$ws.Range("N111").Value = -7766
$ws.Range("H113").Value = 6495
$ws.Range("J113").Value = 6495
$ws.Range("L113").Value = 6495
$ws.Range("N113").Value = -13003
$ws.Range("H118").Value = 166667120
$ws.Range("I118").Value = 200000460
$ws.Range("J118").Value = 420
$ws.Range("K118").Value = 600001380
$ws.Range("L118").Value = 1260
$ws.Range("M118").Value = -599999723
$ws.Range("N118").Value = -4574
$ws.Range("H135").Value = 1411.3125
$ws.Range("I135").Value = 1053
$ws.Range("K135").Value = 9477
$ws.Range("M135").Value = -6942
$ws.Range("H137").Value = 181958.3
$ws.Range("I137").Value = 359335.2
$ws.Range("K137").Value = 1078005.6
$ws.Range("M137").Value = -1075455.6
$ws.Range("H138").Value = 2928.7917
$ws.Range("J138").Value = 3241.3518
$ws.Range("L138").Value = 9724.055399999999
$ws.Range("N138").Value = -20004.0554

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 2999.5
$ws.Range("I8").Value = 1000
$ws.Range("J8").Value = 4999
$ws.Range("K8").Value = 1000
$ws.Range("L8").Value = 4999
$ws.Range("N8").Value = -5287
$ws.Range("M8").Value = -856
$ws.Range("H32").Value = 4684.7314
$ws.Range("I32").Value = 2549.3035
$ws.Range("K32").Value = 2549.3035
$ws.Range("M32").Value = -2262.3035
$ws.Range("H45").Value = 81094.16
$ws.Range("J45").Value = 6213.2856
$ws.Range("L45").Value = 6213.2856
$ws.Range("N45").Value = -6967.2856
$ws.Range("H74").Value = 51056.875
$ws.Range("I74").Value = 6889.968
$ws.Range("K74").Value = 6889.968
$ws.Range("M74").Value = -6015.968
$ws.Range("H77").Value = 51056.875
$ws.Range("I77").Value = 6889.968
$ws.Range("K77").Value = 34449.84
$ws.Range("M77").Value = -30081.84
$ws.Range("H96").Value = 60183.332
$ws.Range("J96").Value = 60183.332
$ws.Range("L96").Value = 60183.332
$ws.Range("N96").Value = -65675.33199999999
$ws.Range("H122").Value = 4692.2856
$ws.Range("I122").Value = 3483.6667
$ws.Range("K122").Value = 10451.0001
$ws.Range("M122").Value = -8001.000100000001
$ws.Range("H132").Value = 2407.1875
$ws.Range("I132").Value = 1793.3334
$ws.Range("K132").Value = 5380.0002
$ws.Range("M132").Value = -2850.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2251.5
$ws.Range("I20").Value = 1653.1666
$ws.Range("J20").Value = 2849.8333
$ws.Range("K20").Value = 1653.1666
$ws.Range("L20").Value = 2849.8333
$ws.Range("M20").Value = -1406.1666
$ws.Range("N20").Value = -3343.8333
$ws.Range("H94").Value = 15447
$ws.Range("I94").Value = 2394
$ws.Range("K94").Value = 2394
$ws.Range("M94").Value = -1943
$ws.Range("H105").Value = 3050.8
$ws.Range("I105").Value = 3050.8
$ws.Range("K105").Value = 3050.8
$ws.Range("M105").Value = -1303.8
$ws.Range("H134").Value = 5532.524
$ws.Range("I134").Value = 2726.5454
$ws.Range("K134").Value = 8179.6362
$ws.Range("M134").Value = -5644.6362

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 500
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 500
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 500
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -840
$ws.Range("H16").Value = 2154.5454
$ws.Range("I16").Value = 1726.25
$ws.Range("K16").Value = 1726.25
$ws.Range("M16").Value = -1439.25
$ws.Range("H17").Value = 598.8
$ws.Range("I17").Value = 598.8
$ws.Range("K17").Value = 598.8
$ws.Range("M17").Value = -424.8
$ws.Range("H58").Value = 1994.9166
$ws.Range("I58").Value = 2020
$ws.Range("K58").Value = 2020
$ws.Range("M58").Value = -1817
$ws.Range("H107").Value = 2643.8667
$ws.Range("J107").Value = 2562.25
$ws.Range("L107").Value = 2562.25
$ws.Range("N107").Value = -6402.25
$ws.Range("H113").Value = 2154.5454
$ws.Range("I113").Value = 1726.25
$ws.Range("K113").Value = 1726.25
$ws.Range("M113").Value = 443.75
$ws.Range("H122").Value = 3261.1765
$ws.Range("I122").Value = 2798.75
$ws.Range("J122").Value = 4371
$ws.Range("K122").Value = 8396.25
$ws.Range("L122").Value = 13113
$ws.Range("M122").Value = -5946.25
$ws.Range("N122").Value = -18013
$ws.Range("H136").Value = 1994.9166
$ws.Range("I136").Value = 2020
$ws.Range("K136").Value = 6060
$ws.Range("M136").Value = -3510

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 255.63829
$ws.Range("I2").Value = 129.35294
$ws.Range("J2").Value = 327.2
$ws.Range("K2").Value = 776.1176399999999
$ws.Range("L2").Value = 1963.2
$ws.Range("M2").Value = -663.1176399999999
$ws.Range("N2").Value = -2189.2
$ws.Range("H17").Value = 779.4286
$ws.Range("I17").Value = 363.75
$ws.Range("J17").Value = 1333.6666
$ws.Range("K17").Value = 1091.25
$ws.Range("L17").Value = 4000.9998
$ws.Range("M17").Value = -922.25
$ws.Range("N17").Value = -4338.9998
$ws.Range("H86").Value = 146
$ws.Range("I86").Value = 146
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 438
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 748
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 146
$ws.Range("I89").Value = 146
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 1314
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 4614
$ws.Range("N89").ClearContents()
$ws.Range("H128").Value = 188554
$ws.Range("I128").Value = 188554
$ws.Range("K128").Value = 565662
$ws.Range("M128").Value = -560682
$ws.Range("H133").Value = 2352.5557
$ws.Range("I133").Value = 2352.5557
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 7057.6671
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -1997.6671
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 1500
$ws.Range("J23").Value = 1500
$ws.Range("L23").Value = 1500
$ws.Range("N23").Value = -1946
$ws.Range("H70").Value = 9338
$ws.Range("I70").Value = 10227.667
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 10227.667
$ws.Range("L70").Value = 4000
$ws.Range("M70").Value = -9957.666999999999
$ws.Range("N70").Value = -4540
$ws.Range("H73").Value = 9338
$ws.Range("I73").Value = 10227.667
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 10227.667
$ws.Range("L73").Value = 4000
$ws.Range("M73").Value = -9291.666999999999
$ws.Range("N73").Value = -5872
$ws.Range("H102").Value = 115753.336
$ws.Range("I102").Value = 3598.3333
$ws.Range("J102").Value = 171830.83
$ws.Range("K102").Value = 3598.3333
$ws.Range("L102").Value = 171830.83
$ws.Range("M102").Value = -1976.3333
$ws.Range("N102").Value = -175074.83
$ws.Range("H122").Value = 1488362.5
$ws.Range("I122").Value = 2967725.2
$ws.Range("J122").Value = 8999.666999999999
$ws.Range("K122").Value = 8903175.600000001
$ws.Range("L122").Value = 26999.001
$ws.Range("M122").Value = -8900725.600000001
$ws.Range("N122").Value = -31899.001
$ws.Range("H126").Value = 4200
$ws.Range("J126").Value = 4830.6
$ws.Range("L126").Value = 14491.8
$ws.Range("N126").Value = -19431.8
$ws.Range("H132").Value = 4445.5
$ws.Range("I132").Value = 3414.1667
$ws.Range("K132").Value = 10242.5001
$ws.Range("M132").Value = -7712.500100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1220
$ws.Range("I16").Value = 1278
$ws.Range("J16").Value = 1147.5
$ws.Range("K16").Value = 1278
$ws.Range("L16").Value = 1147.5
$ws.Range("M16").Value = -1108
$ws.Range("N16").Value = -1487.5
$ws.Range("H22").Value = 251138.75
$ws.Range("I22").Value = 332185
$ws.Range("J22").Value = 8000
$ws.Range("K22").Value = 332185
$ws.Range("L22").Value = 8000
$ws.Range("M22").Value = -331890
$ws.Range("N22").Value = -8590
$ws.Range("H27").Value = 251138.75
$ws.Range("I27").Value = 332185
$ws.Range("J27").Value = 8000
$ws.Range("K27").Value = 332185
$ws.Range("L27").Value = 8000
$ws.Range("M27").Value = -332078
$ws.Range("N27").Value = -8214
$ws.Range("H46").Value = 6439
$ws.Range("I46").Value = 5166.3335
$ws.Range("J46").Value = 7075.3335
$ws.Range("K46").Value = 5166.3335
$ws.Range("L46").Value = 7075.3335
$ws.Range("M46").Value = -4978.3335
$ws.Range("N46").Value = -7451.3335
$ws.Range("H61").Value = 3246.0527
$ws.Range("I61").Value = 3626
$ws.Range("K61").Value = 3626
$ws.Range("M61").Value = -3424
$ws.Range("H68").Value = 2869.1428
$ws.Range("I68").Value = 1673.5
$ws.Range("J68").Value = 4463.3335
$ws.Range("K68").Value = 1673.5
$ws.Range("L68").Value = 4463.3335
$ws.Range("M68").Value = -924.5
$ws.Range("N68").Value = -5961.3335
$ws.Range("H71").Value = 2869.1428
$ws.Range("I71").Value = 1673.5
$ws.Range("J71").Value = 4463.3335
$ws.Range("K71").Value = 8367.5
$ws.Range("L71").Value = 22316.6675
$ws.Range("M71").Value = -4623.5
$ws.Range("N71").Value = -29804.6675
$ws.Range("H82").Value = 959.8261
$ws.Range("I82").Value = 817.1
$ws.Range("J82").Value = 1069.6154
$ws.Range("K82").Value = 817.1
$ws.Range("L82").Value = 1069.6154
$ws.Range("M82").Value = -456.1
$ws.Range("N82").Value = -1791.6154
$ws.Range("H85").Value = 959.8261
$ws.Range("I85").Value = 817.1
$ws.Range("J85").Value = 1069.6154
$ws.Range("K85").Value = 817.1
$ws.Range("L85").Value = 1069.6154
$ws.Range("M85").Value = 430.9
$ws.Range("N85").Value = -3565.6154
$ws.Range("H113").Value = 3246.0527
$ws.Range("I113").Value = 3626
$ws.Range("K113").Value = 3626
$ws.Range("M113").Value = -1456
$ws.Range("H122").Value = 8681.857
$ws.Range("I122").Value = 4995
$ws.Range("K122").Value = 14985
$ws.Range("M122").Value = -12535
$ws.Range("H132").Value = 7114.4136
$ws.Range("I132").Value = 9215.786
$ws.Range("K132").Value = 27647.358
$ws.Range("M132").Value = -25117.358
$ws.Range("H136").Value = 60258.445
$ws.Range("I136").Value = 87147.086
$ws.Range("J136").Value = 6481.1665
$ws.Range("K136").Value = 261441.258
$ws.Range("L136").Value = 19443.4995
$ws.Range("M136").Value = -258891.258
$ws.Range("N136").Value = -24543.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2670.2856
$ws.Range("I122").Value = 1740.4
$ws.Range("K122").Value = 5221.200000000001
$ws.Range("M122").Value = -2771.200000000001
$ws.Range("H132").Value = 128723.125
$ws.Range("I132").Value = 38001
$ws.Range("J132").Value = 183156.4
$ws.Range("K132").Value = 114003
$ws.Range("L132").Value = 549469.2
$ws.Range("M132").Value = -111473
$ws.Range("N132").Value = -554529.2
$ws.Range("H136").Value = 5273.75
$ws.Range("J136").Value = 5384.7144
$ws.Range("L136").Value = 16154.1432
$ws.Range("N136").Value = -21254.1432

